# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Tue Aug 20 23:56:56 UTC 2024 with GitHub Actions".
# Numeric-looking text in column D (Price) is written with a leading "'"
# quote-prefix so Excel keeps storing it as literal text (matching the
# workbook's existing inlineStr cells) instead of silently coercing it to a
# number and dropping significant trailing zeros (e.g. "4.00" -> 4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "58.997.63"
$ws.Range("E2").Value = "  -0.48%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.570.88"
$ws.Range("E3").Value = "  -1.79%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5: BNB
$ws.Range("D5").Value = "'569.33"
$ws.Range("E5").Value = "  +1.75%  "

# Row 6: Solana
$ws.Range("D6").Value = "'142.47"
$ws.Range("E6").Value = "  -1.26%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.24%  "

# Row 8: XRP
$ws.Range("D8").Value = "'0.595"
$ws.Range("E8").Value = "  -0.49%  "

# Row 9: LidoStakedEther
$ws.Range("D9").Value = "2.576.25"
$ws.Range("E9").Value = "  -2.28%  "

# Row 10: Toncoin
$ws.Range("D10").Value = "'6.67"
$ws.Range("E10").Value = "  -2.22%  "

# Row 11: Dogecoin
$ws.Range("E11").Value = "  +1.72%  "

# Row 12: TRON
$ws.Range("D12").Value = "'0.161"
$ws.Range("E12").Value = "  +12.53%  "

# Row 13: Cardano
$ws.Range("E13").Value = "  +2.10%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.020.91"
$ws.Range("E14").Value = "  -1.81%  "

# Row 15: WrappedBTC
$ws.Range("D15").Value = "59.054.90"
$ws.Range("E15").Value = "  -0.27%  "

# Row 16: Avalanche
$ws.Range("D16").Value = "'22.28"
$ws.Range("E16").Value = "  +5.01%  "

# Row 17: ShibaInu
$ws.Range("E17").Value = "  +2.44%  "

# Row 18: WrappedEther
$ws.Range("D18").Value = "2.575.13"
$ws.Range("E18").Value = "  -1.80%  "

# Row 19: Polkadot
$ws.Range("E19").Value = "  +1.06%  "

# Row 20: BitcoinCash
$ws.Range("D20").Value = "'335.96"
$ws.Range("E20").Value = "  -0.80%  "

# Row 21: Chainlink
$ws.Range("D21").Value = "'10.24"
$ws.Range("E21").Value = "  +0.62%  "

# Row 22: Uniswap
$ws.Range("D22").Value = "'6.27"
$ws.Range("E22").Value = "  +1.12%  "

# Row 23: Dai
$ws.Range("E23").Value = "  +0.11%  "

# Row 24: Litecoin
$ws.Range("D24").Value = "'64.46"
$ws.Range("E24").Value = "  -2.61%  "

# Row 25: Polygon
$ws.Range("E25").Value = "  +7.48%  "

# Row 26: Binance-PegBSC-USD
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.63%  "

# Row 27: Kaspa
$ws.Range("E27").Value = "  -3.02%  "

# Row 28: InternetComputer(DFINITY)
$ws.Range("D28").Value = "'7.27"
$ws.Range("E28").Value = "  +0.67%  "

# Row 29: PEPE
$ws.Range("D29").Value = "0.0₃0778"
$ws.Range("E29").Value = "  +1.06%  "

# Row 30: USDe
$ws.Range("E30").Value = "  +0.08%  "

# Row 31: PancakeSwap
$ws.Range("E31").Value = "  -0.41%  "

# Row 32: Monero
$ws.Range("D32").Value = "'159.28"
$ws.Range("E32").Value = "  +2.82%  "

# Row 33: Aptos
$ws.Range("D33").Value = "'6.06"
$ws.Range("E33").Value = "  +0.31%  "

# Row 34: EthereumClassic
$ws.Range("D34").Value = "'18.94"
$ws.Range("E34").Value = "  -0.09%  "

# Row 35: NEARProtocol
$ws.Range("D35").Value = "'4.00"
$ws.Range("E35").Value = "  +0.11%  "

# Row 36: ImmutableX
$ws.Range("D36").Value = "'1.16"
$ws.Range("E36").Value = "  +1.76%  "

# Row 37: SuiNetwork
$ws.Range("D37").Value = "'0.872"
$ws.Range("E37").Value = "  -3.63%  "

# Row 38: Fetch.AI
$ws.Range("E38").Value = "  -4.62%  "

# Row 39: OKB
$ws.Range("D39").Value = "'37.40"
$ws.Range("E39").Value = "  +0.49%  "

# Row 40: Stacks
$ws.Range("E40").Value = "  +0.79%  "

# Row 41: Bittensor
$ws.Range("D41").Value = "'294.64"
$ws.Range("E41").Value = "  +2.90%  "

# Row 42: Filecoin
$ws.Range("E42").Value = "  +0.90%  "

# Row 43: FirstDigitalUSD
$ws.Range("E43").Value = "  +0.32%  "

# Row 44: Aave
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "'0.0978"
$ws.Range("E44").Value = "  +1.86%  "

# Row 45: Stellar
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'129.88"
$ws.Range("E45").Value = "  +10.16%  "

# Row 46: Mantle
$ws.Range("D46").Value = "'0.592"
$ws.Range("E46").Value = "  -1.83%  "

# Row 47: Hedera
$ws.Range("E47").Value = "  -1.17%  "

# Row 48: WhiteBITCoin
$ws.Range("D48").Value = "'10.63"
$ws.Range("E48").Value = "  +0.18%  "

# Row 49: EnergySwap
$ws.Range("D49").Value = "'19.13"
$ws.Range("E49").Value = "  +1.01%  "

# Row 50: VeChain
$ws.Range("E50").Value = "  +1.85%  "

# Row 51: Maker
$ws.Range("D51").Value = "1.942.05"
$ws.Range("E51").Value = "  -0.76%  "
